$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Draw Select tile" column D
$ws.Range("D1").Value = "Horizontal เริ่มที่ค่า"

$ws.Range("D2").Value = 2147483648

# [Fixed] Flip & Rotate tile text fix
$ws.Range("A11").Value = "270 Degree Flip Vertical"

# Update selected cell
$ws.Range("B9").Select()
